# Generate Report for Handoff
# The "aa466399-17a0-42f0-930e-ad0c25aef0fb" entry moves from
# "Handed back: in sync with en-US" to "Ready for handoff", and the
# latest-handoff timestamps for each locale are refreshed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-09 20:44:47"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-09 20:44:52"
